# Adapt output directory style as maven and free.
#
# Fills in the previously-empty rows 24-28 of the "constants" sheet
# (definition list rows, No.5-No.8) with four new constant definitions:
#   TARGET_STYLE_BLANCO / TARGET_STYLE_MAVEN / TARGET_STYLE_FREE
#   TARGET_DIR_SUFFIX_BLANCO / TARGET_DIR_SUFFIX_MAVEN

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 28 currently carries the "bottom border" row style (matching rows
# 29/30, the still-empty trailing rows). Once it gets real data it should
# look like the other populated rows (24-27), so copy that formatting over
# for columns A:E before writing the values (column F is left untouched).
$ws.Range("A24:E24").Copy($ws.Range("A28:E28"))

# No.5
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "TARGET_STYLE_BLANCO"
$ws.Range("C24").Value = "java.lang.String"
$ws.Range("D24").Value = """blanco"""
$ws.Range("E24").Value = "targetdirに設定される文字列"

# No.6
$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "TARGET_STYLE_MAVEN"
$ws.Range("C25").Value = "java.lang.String"
$ws.Range("D25").Value = """maven"""
$ws.Range("E25").Value = "targetdirに設定される文字列"

# No.7
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "TARGET_STYLE_FREE"
$ws.Range("C26").Value = "java.lang.String"
$ws.Range("D26").Value = """free"""
$ws.Range("E26").Value = "targetdirに設定される文字列"

# No.8
$ws.Range("A27").Value = 8
$ws.Range("B27").Value = "TARGET_DIR_SUFFIX_BLANCO"
$ws.Range("C27").Value = "java.lang.String"
$ws.Range("D27").Value = """main"""
$ws.Range("E27").Value = "生成したソースコードを保管するディレクトリのsuffix"

# No.8 (sic - same number reused in the source workbook)
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "TARGET_DIR_SUFFIX_MAVEN"
$ws.Range("C28").Value = "java.lang.String"
$ws.Range("D28").Value = """main/java"""
$ws.Range("E28").Value = "生成したソースコードを保管するディレクトリのsuffix"

# Leave the selection where the author ended up after this edit.
$ws.Range("E29").Select() | Out-Null
